$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the "categoria" column range (C2:C6) and set every cell to "Gelish"
$range = $ws.Range("C2:C6")
$range.Value = "Gelish"
$range.Select()
